$wb = $excel.ActiveWorkbook

# Rename the two dFBA sheets to match the wc_lang rename:
#   DfbaNetReaction -> DfbaObjReaction, DfbaNetSpecies -> DfbaObjSpecies
$wb.Worksheets.Item("dFBA net reactions").Name = "dFBA objective reactions"
$wb.Worksheets.Item("dFBA net species").Name = "dFBA objective species"

# The "dFBA objective species" sheet has a column header that echoed the old name -
# update its text to match the rename.
$wsObjSpecies = $wb.Worksheets.Item("dFBA objective species")
$wsObjSpecies.Range("C1").Value = "dFBA objective reaction"

# Make the renamed "dFBA objective species" sheet the active sheet/selection
# (previously "Environment" was the active sheet).
$wsObjSpecies.Activate() | Out-Null
$wsObjSpecies.Range("E6").Select() | Out-Null
